# The workbook's sheets were renamed, replacing spaces with underscores:
#   "Figura 1"   -> "Figura_1"
#   "Figura 2 a" -> "Figura_2_a"
#   "Figura 2 b" -> "Figura_2_b"
#   "Figura 3"   -> "Figura_3"
#   "Figura 4"   -> "Figura_4"
# (the corresponding change in the selection/activeCellId attributes on
# each sheet is a side-effect of the save round-trip and needs no extra
# work, except on the last sheet where the selected range is collapsed
# from "V56:V58" down to just "V56").

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "Figura_1"
$wb.Worksheets.Item(2).Name = "Figura_2_a"
$wb.Worksheets.Item(3).Name = "Figura_2_b"
$wb.Worksheets.Item(4).Name = "Figura_3"
$wb.Worksheets.Item(5).Name = "Figura_4"

# "Figura 4" (now "Figura_4") was the active sheet with the selection
# anchored at V56 but spanning V56:V58; narrow the selection down to the
# single cell V56.
$ws = $wb.Worksheets.Item("Figura_4")
$ws.Activate()
$ws.Range("V56").Select()
